# Apply updates to Sheet1 of the workbook as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: W(kg) and Numb updated
$ws.Range("G7").Value = 5.057
$ws.Range("H7").Value = 44

# Row 14: W(kg) and Numb updated
$ws.Range("G14").Value = 5.041
$ws.Range("H14").Value = 43

# Rows 22, 24, 25, 36, 37, 38, 39: Numb (H) changed from 0 to -1
$rowsWithNumbChange = @(22, 24, 25, 36, 37, 38, 39)
foreach ($r in $rowsWithNumbChange) {
    $ws.Cells.Item($r, 8).Value = -1
}

# Rows 15-40: RF (I) changed from 7.920823529411765 to 7.15635294117647
for ($r = 15; $r -le 40; $r++) {
    $ws.Cells.Item($r, 9).Value = 7.15635294117647
}
